$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.973.40'
$ws.Range("E2").Value = '  -2.33%  '

$ws.Range("D3").Value = '3.488.48'
$ws.Range("E3").Value = '  +1.30%  '

$ws.Range("E4").Value = '  +0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '582.80'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.44%  '

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '172.92'
$c.ClearFormats()
$ws.Range("E6").Value = '  -3.45%  '

$ws.Range("E7").Value = '  +0.03%  '

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.598'
$c.ClearFormats()
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").Value = '3.488.37'
$ws.Range("E9").Value = '  +1.38%  '

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.131'
$c.ClearFormats()
$ws.Range("E10").Value = '  -5.38%  '

$ws.Range("E11").Value = '  -1.66%  '

$ws.Range("E12").Value = '  -3.65%  '

$ws.Range("D13").Value = '4.093.10'
$ws.Range("E13").Value = '  +1.24%  '

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '29.86'
$c.ClearFormats()
$ws.Range("E15").Value = '  -6.45%  '

$ws.Range("D16").Value = '66.025.62'
$ws.Range("E16").Value = '  -2.20%  '

$ws.Range("E17").Value = '  -2.85%  '

$ws.Range("D18").Value = '3.490.08'
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("E19").Value = '  -3.07%  '

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '13.92'
$c.ClearFormats()
$ws.Range("E20").Value = '  -0.21%  '

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '366.76'
$c.ClearFormats()
$ws.Range("E21").Value = '  -4.87%  '

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '7.73'
$c.ClearFormats()
$ws.Range("E22").Value = '  -0.99%  '

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '72.85'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.31%  '

$ws.Range("E24").Value = '  +0.21%  '

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '0.0000126'
$c.ClearFormats()
$ws.Range("E25").Value = '  +6.07%  '

$ws.Range("E26").Value = '  +0.62%  '

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '9.62'
$c.ClearFormats()
$ws.Range("E27").Value = '  -5.54%  '

$ws.Range("E28").Value = '  +2.53%  '

$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E29").Value = '  +0.10%  '

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '24.13'
$c.ClearFormats()
$ws.Range("E30").Value = '  +2.96%  '

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '5.77'
$c.ClearFormats()
$ws.Range("E31").Value = '  -4.37%  '

$ws.Range("E32").Value = '  -2.76%  '

$ws.Range("E34").Value = '  -0.72%  '

$ws.Range("E35").Value = '  -5.73%  '

$ws.Range("E36").Value = '  -1.14%  '

$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '29.77'
$c.ClearFormats()
$ws.Range("E37").Value = '  +15.85%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '160.84'
$c.ClearFormats()
$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("E39").Value = '  +1.21%  '

$ws.Range("D40").Value = '2.825.13'

$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '1.76'
$c.ClearFormats()
$ws.Range("E41").Value = '  -4.88%  '

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '6.49'
$c.ClearFormats()
$ws.Range("E42").Value = '  -1.66%  '

$ws.Range("E43").Value = '  -6.62%  '

$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '4.46'
$c.ClearFormats()
$ws.Range("E44").Value = '  -1.37%  '

$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '0.0683'
$c.ClearFormats()
$ws.Range("E45").Value = '  -3.73%  '

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '40.08'
$c.ClearFormats()
$ws.Range("E46").Value = '  -2.66%  '

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '24.13'
$c.ClearFormats()
$ws.Range("E47").Value = '  -6.70%  '

$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '0.0288'
$c.ClearFormats()
$ws.Range("E48").Value = '  -2.60%  '

$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '325.27'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("E50").Value = '  -2.20%  '

$ws.Range("E51").Value = '  -2.41%  '
